$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename "Teste" -> "TesteComSucesso" ---------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TesteComSucesso"

# Values (order matters for shared-string table: EMAIL, SENHA, email, teste123)
$ws1.Range("A1").Value = "EMAIL"
$ws1.Range("B1").Value = "SENHA"
$ws1.Range("A2").Value = "teste1@teste.com"
$ws1.Range("B2").Value = "teste123"

# Column widths
$ws1.Columns.Item(1).ColumnWidth = 25.16
$ws1.Columns.Item(2).ColumnWidth = 26.16

# Header row formatting (bold white font on themed blue fill, centered, bordered)
$hdr1 = $ws1.Range("A1:B1")
$hdr1.HorizontalAlignment = -4108
$hdr1.VerticalAlignment = -4108
$f1 = $hdr1.Font
$f1.Bold = $true
$f1.ThemeColor = 2
$hdr1.Interior.ThemeColor = 9

# Data row formatting (bordered + centered)
$row2a = $ws1.Range("A2")
$row2a.HorizontalAlignment = -4108
$row2a.VerticalAlignment = -4108
$row2a.Interior.Pattern = -4142

$row2b = $ws1.Range("B2")
$row2b.HorizontalAlignment = -4108
$row2b.VerticalAlignment = -4108

# Selection for this sheet (it will no longer be the active tab)
$ws1.Range("B7").Select()

# --- Sheet 2: new sheet "TesteComFalha", placed right after sheet 1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "TesteComFalha"

$ws2.Range("A1").Value = "EMAIL"
$ws2.Range("B1").Value = "SENHA"
$ws2.Range("A2").Value = "teste1@teste.com"
$ws2.Range("B2").Value = "teste1234"

$ws2.Columns.Item(1).ColumnWidth = 25.16
$ws2.Columns.Item(2).ColumnWidth = 26.16

$hdr2 = $ws2.Range("A1:B1")
$hdr2.HorizontalAlignment = -4108
$hdr2.VerticalAlignment = -4108
$f2 = $hdr2.Font
$f2.Bold = $true
$f2.ThemeColor = 2
$hdr2.Interior.ThemeColor = 9

$row2a2 = $ws2.Range("A2")
$row2a2.HorizontalAlignment = -4108
$row2a2.VerticalAlignment = -4108
$row2a2.Interior.Pattern = -4142

$row2b2 = $ws2.Range("B2")
$row2b2.HorizontalAlignment = -4108
$row2b2.VerticalAlignment = -4108

# This is the sheet that ends up active / selected
$ws2.Range("B2").Select()
